$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy formatting from column E (the shifted original column D) into new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D with the new fiscal-year data
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 4178000
$ws.Range("D9").Value = 3524000
$ws.Range("D10").Value = 654000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 41000
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 3869000
$ws.Range("D18").Value = 309000
$ws.Range("D20").Value = 31000
$ws.Range("D21").Value = 424000
$ws.Range("D22").Value = 62000
$ws.Range("D23").Value = 278000
$ws.Range("D24").Value = 60000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 218000
$ws.Range("D27").Value = 209000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -92000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -31000
$ws.Range("D33").Value = 117000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 117000
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 115000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 604000
$ws.Range("D44").Value = 477000
$ws.Range("D45").Value = 30000
$ws.Range("D46").Value = 1226000
$ws.Range("D47").Value = 178000
$ws.Range("D48").Value = 483000
$ws.Range("D49").Value = 465000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 374000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2726000
$ws.Range("D57").Value = 700000
$ws.Range("D58").Value = 94000
$ws.Range("D59").Value = 290000
$ws.Range("D60").Value = 1084000
$ws.Range("D61").Value = 730000
$ws.Range("D62").Value = 595000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2439000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 200000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 287000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 117000
$ws.Range("D83").Value = 84000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 251000
$ws.Range("D91").Value = -104000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 111000
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -329000
$ws.Range("D101").Value = -6000
$ws.Range("D102").Value = 27000
